$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the date value on row 77 (was a stray intraday fraction, should align
# with the rest of the date column at 07:00 = .2916666667)
$ws.Range("A77").Value = 45448.2916666667

# Append the new data row (row 78) coming from the latest R script run
$ws.Range("A78").Font.Name = $ws.Range("A77").Font.Name
$ws.Range("A78").Font.Size = $ws.Range("A77").Font.Size
$ws.Range("A78").Font.Color = $ws.Range("A77").Font.Color
$ws.Range("A78").NumberFormat = $ws.Range("A77").NumberFormat
$ws.Range("A78").Value = 45449.332962963

$ws.Range("B78").Value = 900
$ws.Range("C78").Value = 6.30000019073486
$ws.Range("D78").Value = 6.26000022888184
$ws.Range("E78").Value = 6.30000019073486
$ws.Range("F78").Value = 6.26000022888184

# adj_close / ticker columns are stored as text in this sheet (quirk of the
# R export), so force the numeric-looking adj_close into text too, then
# drop the auto-applied quote-prefix formatting so the cell keeps the same
# (default) style as the rest of column G.
$ws.Range("G78").Value = "'6.26000022888184"
$ws.Range("G78").ClearFormats()
$ws.Range("H78").Value = "PAL.MI"
